$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 12:32:16"
$wsZhCn.Range("H4").Value = "2016-03-12 12:32:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 12:32:19"
$wsDeDe.Range("H4").Value = "2016-03-12 12:32:39"
